# daily auto push: 2026-02-01 05:09 UTC
#
# A new sample row for 2026/02/01 (時刻=13, ランキング=149) is inserted
# right after the existing two 2026/02/01 rows (at sheet row 761),
# pushing the previously-logged rows 761-802 down to 762-803.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 761..802 down to 762..803 by inserting a blank row at 761.
$ws.Rows.Item(761).Insert()

# Force column A to be treated as text so the "yyyy/mm/dd"-looking date
# string is not auto-converted into a date serial number, then restore
# the default (Normal) style so the cell doesn't pick up a stray format.
$ws.Cells.Item(761, 1).NumberFormat = "@"
$ws.Cells.Item(761, 1).Value = "2026/02/01"
$ws.Cells.Item(761, 1).Style = "Normal"

$ws.Cells.Item(761, 2).Value = "日"
$ws.Cells.Item(761, 3).Value = 13
$ws.Cells.Item(761, 4).Value = 149
